# "update scripts wuth new tpm"
#
# The NATMI Cort->Sstr4 sheet was regenerated with updated TPM input data.
# The sending/receiving cluster labels, ligand/receptor symbols and the
# categorical columns (A:F) are unchanged; only the derived numeric
# columns (G:T, i.e. expression values, detection/specificity scores and
# edge weights) were recomputed with the new TPM numbers. This script
# writes the recalculated values into the existing cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FAPs -> Cort/Sstr4 -> ECs)
$ws.Range("I2").Value = 0.9158046909573684
$ws.Range("J2").Value = 0.9158046909573684
$ws.Range("M2").Value = 5.415507666666667
$ws.Range("N2").Value = 16.246523
$ws.Range("O2").Value = 0.8321390904960287
$ws.Range("P2").Value = 0.8321390904960289
$ws.Range("Q2").Value = 1.400865471521111
$ws.Range("R2").Value = 12.60778924369
$ws.Range("S2").Value = 0.7620768826052612
$ws.Range("T2").Value = 0.7620768826052613

# Row 3 (FAPs -> Cort/Sstr4 -> FAPs)
$ws.Range("I3").Value = 0.9158046909573684
$ws.Range("J3").Value = 0.9158046909573684
$ws.Range("O3").Value = 0.0334775384739257
$ws.Range("P3").Value = 0.03347753847392571
$ws.Range("S3").Value = 0.03065888677612693
$ws.Range("T3").Value = 0.03065888677612694

# Row 4 (FAPs -> Cort/Sstr4 -> Inflammatory-Mac)
$ws.Range("I4").Value = 0.9158046909573684
$ws.Range("J4").Value = 0.9158046909573684
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.2013576666666667
$ws.Range("N4").Value = 0.6040730000000001
$ws.Range("O4").Value = 0.03094032838984733
$ws.Range("P4").Value = 0.03094032838984733
$ws.Range("Q4").Value = 0.05208653002111112
$ws.Range("R4").Value = 0.4687787701900001
$ws.Range("S4").Value = 0.02833529787918362
$ws.Range("T4").Value = 0.02833529787918363

# Row 5 (FAPs -> Cort/Sstr4 -> MuSCs)
$ws.Range("I5").Value = 0.9158046909573684
$ws.Range("J5").Value = 0.9158046909573684
$ws.Range("M5").Value = 0.200945
$ws.Range("N5").Value = 0.602835
$ws.Range("O5").Value = 0.03087691862555289
$ws.Range("P5").Value = 0.03087691862555289
$ws.Range("Q5").Value = 0.05197978278333334
$ws.Range("R5").Value = 0.46781804505
$ws.Range("S5").Value = 0.02827722691959028
$ws.Range("T5").Value = 0.02827722691959028

# Row 6 (FAPs -> Cort/Sstr4 -> Resolving-Mac)
$ws.Range("I6").Value = 0.9158046909573684
$ws.Range("J6").Value = 0.9158046909573684
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.4722556666666667
$ws.Range("N6").Value = 1.416767
$ws.Range("O6").Value = 0.07256612401464529
$ws.Range("P6").Value = 0.07256612401464529
$ws.Range("Q6").Value = 0.1221615216677778
$ws.Range("R6").Value = 1.09945369501
$ws.Range("S6").Value = 0.0664563967772063
$ws.Range("T6").Value = 0.0664563967772063

# Row 7 (Inflammatory-Mac -> Cort/Sstr4 -> ECs)
$ws.Range("G7").Value = 0.02378166666666667
$ws.Range("H7").Value = 0.07134500000000001
$ws.Range("I7").Value = 0.08419530904263167
$ws.Range("J7").Value = 0.08419530904263166
$ws.Range("M7").Value = 5.415507666666667
$ws.Range("N7").Value = 16.246523
$ws.Range("O7").Value = 0.8321390904960287
$ws.Range("P7").Value = 0.8321390904960289
$ws.Range("Q7").Value = 0.1287897981594445
$ws.Range("R7").Value = 1.159108183435
$ws.Range("S7").Value = 0.07006220789076759
$ws.Range("T7").Value = 0.07006220789076759

# Row 8 (Inflammatory-Mac -> Cort/Sstr4 -> FAPs)
$ws.Range("G8").Value = 0.02378166666666667
$ws.Range("H8").Value = 0.07134500000000001
$ws.Range("I8").Value = 0.08419530904263167
$ws.Range("J8").Value = 0.08419530904263166
$ws.Range("O8").Value = 0.0334775384739257
$ws.Range("P8").Value = 0.03347753847392571
$ws.Range("Q8").Value = 0.005181303789444445
$ws.Range("R8").Value = 0.046631734105
$ws.Range("S8").Value = 0.002818651697798766
$ws.Range("T8").Value = 0.002818651697798766

# Row 9 (Inflammatory-Mac -> Cort/Sstr4 -> Inflammatory-Mac)
$ws.Range("G9").Value = 0.02378166666666667
$ws.Range("H9").Value = 0.07134500000000001
$ws.Range("I9").Value = 0.08419530904263167
$ws.Range("J9").Value = 0.08419530904263166
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.2013576666666667
$ws.Range("N9").Value = 0.6040730000000001
$ws.Range("O9").Value = 0.03094032838984733
$ws.Range("P9").Value = 0.03094032838984733
$ws.Range("Q9").Value = 0.004788620909444445
$ws.Range("R9").Value = 0.04309758818500001
$ws.Range("S9").Value = 0.002605030510663706
$ws.Range("T9").Value = 0.002605030510663706

# Row 10 (Inflammatory-Mac -> Cort/Sstr4 -> MuSCs)
$ws.Range("G10").Value = 0.02378166666666667
$ws.Range("H10").Value = 0.07134500000000001
$ws.Range("I10").Value = 0.08419530904263167
$ws.Range("J10").Value = 0.08419530904263166
$ws.Range("M10").Value = 0.200945
$ws.Range("N10").Value = 0.602835
$ws.Range("O10").Value = 0.03087691862555289
$ws.Range("P10").Value = 0.03087691862555289
$ws.Range("Q10").Value = 0.004778807008333334
$ws.Range("R10").Value = 0.04300926307500001
$ws.Range("S10").Value = 0.002599691705962616
$ws.Range("T10").Value = 0.002599691705962616

# Row 11 (Inflammatory-Mac -> Cort/Sstr4 -> Resolving-Mac)
$ws.Range("G11").Value = 0.02378166666666667
$ws.Range("H11").Value = 0.07134500000000001
$ws.Range("I11").Value = 0.08419530904263167
$ws.Range("J11").Value = 0.08419530904263166
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.4722556666666667
$ws.Range("N11").Value = 1.416767
$ws.Range("O11").Value = 0.07256612401464529
$ws.Range("P11").Value = 0.07256612401464529
$ws.Range("Q11").Value = 0.01123102684611111
$ws.Range("R11").Value = 0.101079241615
$ws.Range("S11").Value = 0.006109727237438996
$ws.Range("T11").Value = 0.006109727237438995
